# Update column F ("dSF") values for specific rows per the diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -5
    3  = 4
    11 = -1
    15 = -2
    19 = 4
    26 = 0
    31 = 0
    32 = -1
    34 = -1
    35 = 1
    38 = 0
    42 = 1
    43 = 1
    47 = -1
    49 = 0
    50 = 2
    52 = 1
    61 = 1
    64 = 2
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
